# Add a header (containing the questionnaire number) to the document's
# only section, so the questionnaire number survives printing.
#
# wdHeaderFooterPrimary = 1
# wdAlignParagraphCenter = 1

$d = $word.ActiveDocument

$section = $d.Sections(1)
$header  = $section.Headers(1)

# Writing into the header's Range (rather than assigning .Text, which
# Word treats as "turn on headers/footers for the whole document" and
# mints the full even/first/primary header+footer set) creates just the
# single default header part the section needs.
$header.Range.InsertAfter("Questionnaire 24")

# Style + center the paragraph that now holds the text.
$headerPara = $header.Range.Paragraphs(1)
$headerPara.Range.Style = "Header"
$headerPara.Range.ParagraphFormat.Alignment = 1

# Apply the run-level font formatting to the text itself (not the
# trailing paragraph mark) so it doesn't leak into the paragraph's
# rPr (pilcrow formatting).
$textRange = $header.Range
$textRange.End = $textRange.End - 1
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
